# Rename the inline picture shapes living in the document's headers/footers.
#
# Three logo pictures get their wp:docPr/@name flipped:
#   - Footer (default, section 1)    : "PearsonLogo" picture  image1.png -> image2.png
#   - Footer (first page, section 1) : "PearsonLogo" picture  image1.png -> image2.png
#   - Header (first page, section 1) : "BTec_Logo-Orange"     image2.jpg -> image1.jpg
#
# Setting InlineShape.Name directly on a shape fetched from a
# HeaderFooter.Range sometimes fails to stick for footer stories, so the
# shape is Select()-ed first and then renamed through $word.Selection --
# that path reliably reaches the underlying drawing object.

$d = $word.ActiveDocument

function Rename-LogoInStory($story, $matchAltText, $newName) {
    if ($story -eq $null -or -not $story.Exists) {
        return
    }
    $shapes = $story.Range.InlineShapes
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $candidate = $shapes.Item($k)
        if ($candidate.AlternativeText -eq $matchAltText) {
            $candidate.Select() | Out-Null
            $selected = $word.Selection.InlineShapes.Item(1)
            $selected.Name = $newName
        }
    }
}

foreach ($sec in $d.Sections) {

    # --- Footers: both the default and first-page footer carry the Pearson logo ---
    for ($i = 1; $i -le 3; $i++) {
        Rename-LogoInStory $sec.Footers.Item($i) "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" "image2.png"
    }

    # --- Headers: the first-page header carries the BTec logo ---
    for ($i = 1; $i -le 3; $i++) {
        Rename-LogoInStory $sec.Headers.Item($i) "BTec_Logo-Orange" "image1.jpg"
    }
}
